$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.347.79'
$ws.Range('E2').Value = '  -0.03%  '

$ws.Range('D3').Value = '1.562.01'
$ws.Range('E3').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.007'
$ws.Range('E5').Value = '  +0.49%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.30'
$ws.Range('E6').Value = '  -0.46%  '

$ws.Range('E7').Value = '  +0.62%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.41'
$ws.Range('E8').Value = '  +0.67%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3348'
$ws.Range('E9').Value = '  -1.44%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07374'
$ws.Range('E10').Value = '  -3.40%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.110'
$ws.Range('E11').Value = '  -4.82%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.008'
$ws.Range('E12').Value = '  +0.58%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.60'
$ws.Range('E13').Value = '  -3.74%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.831'
$ws.Range('E14').Value = '  -3.53%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.829'
$ws.Range('E15').Value = '  -1.23%  '

$ws.Range('D16').Value = '1.561.48'
$ws.Range('E16').Value = '  -0.28%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001099'
$ws.Range('E17').Value = '  -2.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '88.72'
$ws.Range('E18').Value = '  -1.43%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06684'
$ws.Range('E19').Value = '  -0.63%  '

$ws.Range('E20').Value = '  +0.48%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.113'
$ws.Range('E21').Value = '  -1.87%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.14'
$ws.Range('E22').Value = '  -2.35%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.77'
$ws.Range('E23').Value = '  -1.96%  '

$ws.Range('D24').Value = '22.356.28'
$ws.Range('E24').Value = '  +0.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.374'
$ws.Range('E25').Value = '  -0.97%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.566'
$ws.Range('E26').Value = '  -9.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.79'
$ws.Range('E27').Value = '  -1.88%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '146.72'
$ws.Range('E28').Value = '  +0.92%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.021'
$ws.Range('E29').Value = '  +0.71%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.20'
$ws.Range('E30').Value = '  -0.91%  '

$ws.Range('D31').Value = '1.733.50'
$ws.Range('E31').Value = '  -0.34%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.018'
$ws.Range('E32').Value = '  +0.52%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9736'
$ws.Range('E33').Value = '  -2.80%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.837'
$ws.Range('E34').Value = '  -5.75%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.685'
$ws.Range('E35').Value = '  -3.45%  '

$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08427'
$ws.Range('E36').Value = '  -0.11%  '

$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.387'
$ws.Range('E37').Value = '  +6.84%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02425'
$ws.Range('E38').Value = '  -4.22%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2243'
$ws.Range('E39').Value = '  -3.42%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06338'
$ws.Range('E40').Value = '  -0.71%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.289'
$ws.Range('E41').Value = '  -4.13%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6155'
$ws.Range('E42').Value = '  -2.98%  '

$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.006'
$ws.Range('E43').Value = '  +0.53%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.83'
$ws.Range('E44').Value = '  -7.33%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.77'
$ws.Range('E45').Value = '  -2.69%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.779'
$ws.Range('E46').Value = '  +0.60%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5712'
$ws.Range('E47').Value = '  -4.30%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.016'
$ws.Range('E48').Value = '  -3.66%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.232'
$ws.Range('E49').Value = '  -2.83%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '123.49'
$ws.Range('E50').Value = '  -0.89%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07288'
$ws.Range('E51').Value = '  +0.30%  '
